# Apply cryptos list price/volume updates (and PEPE/Maker row reorder)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain (unstyled) cell used as a format donor so that cells which must be
# forced back to text after receiving a numeric-looking string keep the
# workbook's original (default) cell style.
$formatDonor = $ws.Range("B2")

$ws.Range("D2").Value = "70.805.77"
$ws.Range("E2").Value = "  +0.04%  "
$ws.Range("D3").Value = "3.532.85"
$ws.Range("E3").Value = "  -1.02%  "
$ws.Range("E4").Value = "  -0.09%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "624.95"
$formatDonor.Copy()
$c.PasteSpecial(-4122)

$ws.Range("E5").Value = "  +2.68%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "175.07"
$formatDonor.Copy()
$c.PasteSpecial(-4122)

$ws.Range("E6").Value = "  +0.55%  "
$ws.Range("D7").Value = "3.530.81"
$ws.Range("E7").Value = "  -1.01%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.610"
$formatDonor.Copy()
$c.PasteSpecial(-4122)

$ws.Range("E8").Value = "  -1.18%  "
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("E10").Value = "  +0.39%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "7.16"
$formatDonor.Copy()
$c.PasteSpecial(-4122)

$ws.Range("E11").Value = "  -5.60%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.587"
$formatDonor.Copy()
$c.PasteSpecial(-4122)

$ws.Range("E12").Value = "  -0.30%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "46.70"
$formatDonor.Copy()
$c.PasteSpecial(-4122)

$ws.Range("E13").Value = "  -0.26%  "
$ws.Range("E14").Value = "  -0.44%  "
$ws.Range("D15").Value = "4.103.11"
$ws.Range("E15").Value = "  -1.17%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "8.43"
$formatDonor.Copy()
$c.PasteSpecial(-4122)

$ws.Range("E16").Value = "  -0.01%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "608.53"
$formatDonor.Copy()
$c.PasteSpecial(-4122)

$ws.Range("E17").Value = "  -1.27%  "
$ws.Range("D18").Value = "3.529.38"
$ws.Range("E18").Value = "  -1.15%  "
$ws.Range("D19").Value = "70.880.67"
$ws.Range("E19").Value = "  +0.06%  "
$ws.Range("E20").Value = "  +1.27%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "17.85"
$formatDonor.Copy()
$c.PasteSpecial(-4122)

$ws.Range("E21").Value = "  +2.44%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "0.887"
$formatDonor.Copy()
$c.PasteSpecial(-4122)

$ws.Range("E22").Value = "  -0.20%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "9.06"
$formatDonor.Copy()
$c.PasteSpecial(-4122)

$ws.Range("E23").Value = "  -4.33%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "15.69"
$formatDonor.Copy()
$c.PasteSpecial(-4122)

$ws.Range("E24").Value = "  -2.59%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "98.21"
$formatDonor.Copy()
$c.PasteSpecial(-4122)

$ws.Range("E25").Value = "  +0.95%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "3.79"
$formatDonor.Copy()
$c.PasteSpecial(-4122)

$ws.Range("E26").Value = "  -0.94%  "
$ws.Range("E27").Value = "  +0.03%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "2.59"
$formatDonor.Copy()
$c.PasteSpecial(-4122)

$ws.Range("E28").Value = "  -1.82%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "33.90"
$formatDonor.Copy()
$c.PasteSpecial(-4122)

$ws.Range("E29").Value = "  +1.02%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "9.13"
$formatDonor.Copy()
$c.PasteSpecial(-4122)

$ws.Range("E30").Value = "  -0.22%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "3.05"
$formatDonor.Copy()
$c.PasteSpecial(-4122)

$ws.Range("E31").Value = "  -0.40%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "8.16"
$formatDonor.Copy()
$c.PasteSpecial(-4122)

$ws.Range("E32").Value = "  -4.17%  "
$ws.Range("E33").Value = "  +0.17%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "6.86"
$formatDonor.Copy()
$c.PasteSpecial(-4122)

$ws.Range("E34").Value = "  -2.59%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "631.51"
$formatDonor.Copy()
$c.PasteSpecial(-4122)

$ws.Range("E35").Value = "  +5.14%  "
$ws.Range("E36").Value = "  -1.73%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "10.86"
$formatDonor.Copy()
$c.PasteSpecial(-4122)

$ws.Range("E37").Value = "  +0.01%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "3.50"
$formatDonor.Copy()
$c.PasteSpecial(-4122)

$ws.Range("E38").Value = "  -5.89%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.0476"
$formatDonor.Copy()
$c.PasteSpecial(-4122)

$ws.Range("E39").Value = "  -1.31%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "56.80"
$formatDonor.Copy()
$c.PasteSpecial(-4122)

$ws.Range("E40").Value = "  -0.97%  "
$ws.Range("E41").Value = "  +0.32%  "
$ws.Range("E42").Value = "  +1.44%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "3.363.74"
$ws.Range("E43").Value = "  -0.81%  "
$ws.Range("B44").Value = "PEPE"
$ws.Range("C44").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D44").Value = "0.0₃0737"
$ws.Range("E44").Value = "  +3.95%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "3.03"
$formatDonor.Copy()
$c.PasteSpecial(-4122)

$ws.Range("E45").Value = "  +0.90%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.313"
$formatDonor.Copy()
$c.PasteSpecial(-4122)

$ws.Range("E46").Value = "  -2.68%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "32.24"
$formatDonor.Copy()
$c.PasteSpecial(-4122)

$ws.Range("E47").Value = "  -3.22%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "2.57"
$formatDonor.Copy()
$c.PasteSpecial(-4122)

$ws.Range("E48").Value = "  -1.95%  "
$ws.Range("E49").Value = "  +0.23%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "133.12"
$formatDonor.Copy()
$c.PasteSpecial(-4122)

$ws.Range("E50").Value = "  -0.31%  "
$ws.Range("E51").Value = "  +5.39%  "
$excel.CutCopyMode = $false
